# updateSettlementBIC.230511.xlsx - Build v2.1.2
# Fix SearchCriteria variants and Schemas sheet grouping/sorting
#
# This script updates the "Body", "200", "204", "400", "401", "403", "404",
# "429" and "500" sheets so that their row 3 (the first data row right below
# the header row) describes the request/response/error schema reference
# used by the API definition, instead of the old "dateTime" / error detail
# rows. Sheets that previously only had the header rows (204, 401, 403,
# 404, 429, 500) gain a new row 3 with the schema reference; sheets that
# previously listed several fields (Body, 200, 400) are collapsed down to
# just the single schema-reference row.

$wb = $excel.ActiveWorkbook

function Set-SchemaRow($SheetName, $SectionValue, $SchemaName) {
    $ws = $wb.Worksheets.Item($SheetName)

    $ws.Range("A3").Value = $SectionValue
    $ws.Range("B3").Value = $SchemaName
    $ws.Range("C3").Value = ""
    $ws.Range("D3").Value = ""
    $ws.Range("E3").Value = "schema"
    $ws.Range("F3").Value = ""
    $ws.Range("G3").Value = $SchemaName
    $ws.Range("H3").Value = ""
    $ws.Range("I3").Value = "Yes"
    $ws.Range("J3").Value = ""
    $ws.Range("K3").Value = ""
    $ws.Range("L3").Value = ""
    $ws.Range("M3").Value = ""
    $ws.Range("N3").Value = ""
    $ws.Range("O3").Value = ""
}

# --- "Body" sheet (Request) ------------------------------------------------
# Drop the old field-by-field rows (4-13) and turn row 3 into the schema
# reference row for the request body. Contents are cleared (rather than
# using a row delete) so that other range references on the sheet (e.g.
# data validations / conditional formatting anchored to A1048576) are not
# shifted.
$wsBody = $wb.Worksheets.Item("Body")
$wsBody.Range("A4:O13").ClearContents()
Set-SchemaRow "Body" "body" "updateSettlementBIC.230511Request"

# --- "200" sheet (Response) -------------------------------------------------
# Drop the old commandRef row (4) and turn row 3 into the schema reference
# row for the success response body.
$ws200 = $wb.Worksheets.Item("200")
$ws200.Range("A4:O4").ClearContents()
Set-SchemaRow "200" "content" "updateSettlementBIC.230511Response"

# --- "204" sheet -------------------------------------------------------------
# No existing data rows; add the new schema reference row.
Set-SchemaRow "204" "content" "updateSettlementBIC.230511Response"

# --- "400" sheet ---------------------------------------------------------------
# Drop the old errorCode/errorCodeDescription/requestId rows (4-6) and turn
# row 3 into the schema reference row for the generic error response.
$ws400 = $wb.Worksheets.Item("400")
$ws400.Range("A4:O6").ClearContents()
Set-SchemaRow "400" "content" "errorResponse"

# --- "401", "403", "404", "429", "500" sheets --------------------------------
# No existing data rows; add the new schema reference row pointing at the
# errorResponse1 schema.
foreach ($sheetName in @("401", "403", "404", "429", "500")) {
    Set-SchemaRow $sheetName "content" "errorResponse1"
}
